$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.660.54'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '1.889.46'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'237.18"
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = "'0.4841"
$ws.Range('E7').Value = '  +1.13%  '
$ws.Range('D8').Value = "'0.2858"
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('D9').Value = "'0.06547"
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').Value = '1.839.61'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').Value = "'0.07444"
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').Value = "'16.56"
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').Value = "'5.088"
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').Value = "'87.85"
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('D15').Value = "'0.6634"
$ws.Range('E15').Value = '  +3.05%  '
$ws.Range('D16').Value = '30.606.02'
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = "'0.000007606"
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').Value = "'229.51"
$ws.Range('E20').Value = '  +2.92%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.080.27'
$ws.Range('E21').Value = '  -0.85%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = "'1.001"
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = "'5.265"
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = "'6.186"
$ws.Range('E24').Value = '  +2.20%  '
$ws.Range('D25').Value = "'9.421"
$ws.Range('E25').Value = '  +2.70%  '
$ws.Range('D26').Value = "'167.58"
$ws.Range('E26').Value = '  +2.61%  '
$ws.Range('D27').Value = "'18.67"
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').Value = "'1.954"
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('D29').Value = "'0.1022"
$ws.Range('E29').Value = '  +11.10%  '
$ws.Range('D30').Value = "'1.395"
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').Value = "'4.331"
$ws.Range('E31').Value = '  +2.07%  '
$ws.Range('D32').Value = "'4.022"
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('E33').Value = '  +2.03%  '
$ws.Range('D34').Value = "'1.203"
$ws.Range('E34').Value = '  +5.39%  '
$ws.Range('D35').Value = "'0.7510"
$ws.Range('E35').Value = '  +3.87%  '
$ws.Range('D36').Value = "'0.9993"
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  +0.89%  '
$ws.Range('D38').Value = "'0.01893"
$ws.Range('E38').Value = '  +3.56%  '
$ws.Range('D39').Value = "'2.667"
$ws.Range('E39').Value = '  +2.79%  '
$ws.Range('E40').Value = '  +2.38%  '
$ws.Range('D41').Value = "'2.059"
$ws.Range('E41').Value = '  +1.17%  '
$ws.Range('D42').Value = "'107.25"
$ws.Range('E42').Value = '  +1.20%  '
$ws.Range('D43').Value = "'0.4266"
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = "'5.636"
$ws.Range('E45').Value = '  -4.40%  '
$ws.Range('D46').Value = "'7.418"
$ws.Range('E46').Value = '  +2.36%  '
$ws.Range('D47').Value = "'64.74"
$ws.Range('E47').Value = '  +2.47%  '
$ws.Range('D48').Value = "'0.1272"
$ws.Range('E48').Value = '  -1.67%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = "'8.938"
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('D51').Value = "'34.03"
$ws.Range('E51').Value = '  +0.97%  '
